# Adds a new "Postgres JSON (disabled 2nd level hibernate cache, lz4 compression)"
# timings row to the "timings" sheet, inserted right after the existing
# "Postgres JSON (disabled 2nd level hibernate cache)" row (row 9) and before
# the blank separator row that used to sit at row 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("timings")

# Insert a brand-new row at row 10 - this pushes the blank spacer row
# (formerly row 10) and everything below it down by one.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row with the new benchmark data.
$ws.Range("A10").Value = "Postgres JSON (disabled 2nd level hibernate cache, lz4 compression)"
$ws.Range("B10").Value = 26
$ws.Range("C10").Value = 271
$ws.Range("D10").Value = 470
$ws.Range("E10").Value = 12
$ws.Range("F10").Value = 130
$ws.Range("G10").Value = 242
$ws.Range("H10").Value = "local docker"

# Column A needs to be widened so the new, longer label fits.
$ws.Columns.Item(1).ColumnWidth = 60.5

# Move / restore the active selection like the source workbook.
$ws.Range("A22").Select()
